# Scheduled runner refresh: update cached market-price figures
# (currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ /
# LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ columns H:N)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 75474.75
$ws.Range("J40").Value = 67299.664
$ws.Range("L40").Value = 67299.664
$ws.Range("N40").Value = -67649.664
$ws.Range("H74").Value = 4624.7
$ws.Range("I74").Value = 4280.875
$ws.Range("J74").Value = 6000
$ws.Range("K74").Value = 4280.875
$ws.Range("L74").Value = 6000
$ws.Range("M74").Value = -3344.875
$ws.Range("N74").Value = -7872
$ws.Range("H77").Value = 4624.7
$ws.Range("I77").Value = 4280.875
$ws.Range("J77").Value = 6000
$ws.Range("K77").Value = 21404.375
$ws.Range("L77").Value = 30000
$ws.Range("M77").Value = -16724.375
$ws.Range("N77").Value = -39360
$ws.Range("H132").Value = 5806.04
$ws.Range("I132").Value = 1983.35
$ws.Range("K132").Value = 5950.049999999999
$ws.Range("M132").Value = -3420.049999999999
$ws.Range("H138").Value = 5412.953
$ws.Range("I138").Value = 1420.9166
$ws.Range("J138").Value = 6983.5903
$ws.Range("K138").Value = 4262.7498
$ws.Range("L138").Value = 20950.7709
$ws.Range("M138").Value = 877.2502000000004
$ws.Range("N138").Value = -31230.7709

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4805.5
$ws.Range("I61").Value = 5050.5557
$ws.Range("K61").Value = 5050.5557
$ws.Range("M61").Value = -4838.5557
$ws.Range("H74").Value = 27780138
$ws.Range("I74").Value = 35716820
$ws.Range("J74").Value = 1750
$ws.Range("K74").Value = 35716820
$ws.Range("L74").Value = 1750
$ws.Range("M74").Value = -35715946
$ws.Range("N74").Value = -3498
$ws.Range("H77").Value = 27780138
$ws.Range("I77").Value = 35716820
$ws.Range("J77").Value = 1750
$ws.Range("K77").Value = 178584100
$ws.Range("L77").Value = 8750
$ws.Range("M77").Value = -178579732
$ws.Range("N77").Value = -17486
$ws.Range("H132").Value = 1830.6428
$ws.Range("I132").Value = 1740.7179
$ws.Range("J132").Value = 2999.6667
$ws.Range("K132").Value = 5222.153700000001
$ws.Range("L132").Value = 8999.000100000001
$ws.Range("M132").Value = -2692.153700000001
$ws.Range("N132").Value = -14059.0001
$ws.Range("H136").Value = 4805.5
$ws.Range("I136").Value = 5050.5557
$ws.Range("K136").Value = 15151.6671
$ws.Range("M136").Value = -12601.6671

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 4027.85
$ws.Range("I105").Value = 3897.111
$ws.Range("K105").Value = 3897.111
$ws.Range("M105").Value = -2150.111

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 957.3333
$ws.Range("J16").Value = 2500
$ws.Range("L16").Value = 2500
$ws.Range("N16").Value = -3074
$ws.Range("H31").Value = 15627418
$ws.Range("I31").Value = 17859172
$ws.Range("J31").Value = 5133.625
$ws.Range("K31").Value = 17859172
$ws.Range("L31").Value = 5133.625
$ws.Range("M31").Value = -17858877
$ws.Range("N31").Value = -5723.625
$ws.Range("H34").Value = 15627418
$ws.Range("I34").Value = 17859172
$ws.Range("J34").Value = 5133.625
$ws.Range("K34").Value = 17859172
$ws.Range("L34").Value = 5133.625
$ws.Range("M34").Value = -17858970
$ws.Range("N34").Value = -5537.625
$ws.Range("H58").Value = 3438
$ws.Range("I58").Value = 2962
$ws.Range("J58").Value = 4072.6667
$ws.Range("K58").Value = 2962
$ws.Range("L58").Value = 4072.6667
$ws.Range("M58").Value = -2759
$ws.Range("N58").Value = -4478.6667
$ws.Range("H109").Value = 46991
$ws.Range("J109").Value = 46991
$ws.Range("L109").Value = 46991
$ws.Range("N109").Value = -49071
$ws.Range("H113").Value = 957.3333
$ws.Range("J113").Value = 2500
$ws.Range("L113").Value = 2500
$ws.Range("N113").Value = -6840
$ws.Range("H132").Value = 148160830
$ws.Range("I132").Value = 190477920
$ws.Range("K132").Value = 571433760
$ws.Range("M132").Value = -571431230
$ws.Range("H134").Value = 2247.465
$ws.Range("I134").Value = 1395.3871
$ws.Range("K134").Value = 4186.1613
$ws.Range("M134").Value = -1651.1613
$ws.Range("H136").Value = 3438
$ws.Range("I136").Value = 2962
$ws.Range("J136").Value = 4072.6667
$ws.Range("K136").Value = 8886
$ws.Range("L136").Value = 12218.0001
$ws.Range("M136").Value = -6336
$ws.Range("N136").Value = -17318.0001
$ws.Range("H141").Value = 123573.62
$ws.Range("J141").Value = 131558.95
$ws.Range("L141").Value = 131558.95
$ws.Range("N141").Value = -141918.95

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 7777827.5
$ws.Range("I4").Value = 125055.75
$ws.Range("K4").Value = 375167.25
$ws.Range("M4").Value = -375055.25
$ws.Range("H5").Value = 497.7857
$ws.Range("I5").Value = 389.25
$ws.Range("J5").Value = 1149
$ws.Range("K5").Value = 1167.75
$ws.Range("L5").Value = 3447
$ws.Range("M5").Value = -1055.75
$ws.Range("N5").Value = -3671
$ws.Range("H113").Value = 741.0909
$ws.Range("I113").Value = 999
$ws.Range("J113").Value = 715.3
$ws.Range("K113").Value = 2997
$ws.Range("L113").Value = 2145.9
$ws.Range("M113").Value = -827
$ws.Range("N113").Value = -6485.9
$ws.Range("H132").Value = 1223.9117
$ws.Range("I132").Value = 1003.96155
$ws.Range("J132").Value = 1938.75
$ws.Range("K132").Value = 9035.65395
$ws.Range("L132").Value = 17448.75
$ws.Range("M132").Value = -6505.65395
$ws.Range("N132").Value = -22508.75
$ws.Range("H135").Value = 497.7857
$ws.Range("I135").Value = 389.25
$ws.Range("J135").Value = 1149
$ws.Range("K135").Value = 3503.25
$ws.Range("L135").Value = 10341
$ws.Range("M135").Value = -968.25
$ws.Range("N135").Value = -15411

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8660.916999999999
$ws.Range("I70").Value = 8550.883
$ws.Range("J70").Value = 8928.143
$ws.Range("K70").Value = 8550.883
$ws.Range("L70").Value = 8928.143
$ws.Range("M70").Value = -8280.883
$ws.Range("N70").Value = -9468.143
$ws.Range("H73").Value = 8660.916999999999
$ws.Range("I73").Value = 8550.883
$ws.Range("J73").Value = 8928.143
$ws.Range("K73").Value = 8550.883
$ws.Range("L73").Value = 8928.143
$ws.Range("M73").Value = -7614.883
$ws.Range("N73").Value = -10800.143
$ws.Range("H80").Value = 132333.33
$ws.Range("I80").Value = 205600.2
$ws.Range("K80").Value = 205600.2
$ws.Range("M80").Value = -204602.2
$ws.Range("H83").Value = 132333.33
$ws.Range("I83").Value = 205600.2
$ws.Range("K83").Value = 1028001
$ws.Range("M83").Value = -1023009
$ws.Range("H122").Value = 253185.75
$ws.Range("J122").Value = 7415.9165
$ws.Range("L122").Value = 22247.7495
$ws.Range("N122").Value = -27147.7495
$ws.Range("H123").Value = 34629.668
$ws.Range("J123").Value = 34629.668
$ws.Range("L123").Value = 34629.668
$ws.Range("N123").Value = -39529.668
$ws.Range("H126").Value = 2296.9167
$ws.Range("I126").Value = 2012
$ws.Range("J126").Value = 3379.6
$ws.Range("K126").Value = 6036
$ws.Range("L126").Value = 10138.8
$ws.Range("M126").Value = -3566
$ws.Range("N126").Value = -15078.8
$ws.Range("H132").Value = 136202
$ws.Range("I132").Value = 223370
$ws.Range("K132").Value = 670110
$ws.Range("M132").Value = -667580
$ws.Range("H136").Value = 29239.533
$ws.Range("J136").Value = 29239.533
$ws.Range("L136").Value = 87718.599
$ws.Range("N136").Value = -92818.599
$ws.Range("H141").Value = 32697.5
$ws.Range("J141").Value = 32697.5
$ws.Range("L141").Value = 32697.5
$ws.Range("N141").Value = -43057.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1351
$ws.Range("I16").Value = 1351
$ws.Range("K16").Value = 1351
$ws.Range("M16").Value = -1181
$ws.Range("H139").Value = 85715
$ws.Range("J139").Value = 85715
$ws.Range("L139").Value = 85715
$ws.Range("N139").Value = -95995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1893.8334
$ws.Range("I96").Value = 2131.6667
$ws.Range("J96").Value = 1656
$ws.Range("K96").Value = 2131.6667
$ws.Range("L96").Value = 1656
$ws.Range("M96").Value = -758.6667000000002
$ws.Range("N96").Value = -4402
$ws.Range("H107").Value = 865.53845
$ws.Range("I107").Value = 850.5
$ws.Range("J107").Value = 872.2222
$ws.Range("K107").Value = 2551.5
$ws.Range("L107").Value = 2616.6666
$ws.Range("M107").Value = -631.5
$ws.Range("N107").Value = -6456.6666
$ws.Range("H123").Value = 84294.5
$ws.Range("J123").Value = 84294.5
$ws.Range("L123").Value = 84294.5
$ws.Range("N123").Value = -94094.5
$ws.Range("H140").Value = 95947.336
$ws.Range("J140").Value = 95947.336
$ws.Range("L140").Value = 95947.336
$ws.Range("N140").Value = -106307.336
